$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for engine type columns M:S
$headers = @("dohcv", "ohcv", "ohc", "l", "rotor", "ohcf", "dohc")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 13 + $i   # M=13 ... S=19
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# New data values for rows 2-5, columns M:S
$data = @(
    @(3462, 15841, 17602, 14139, 0, 1196, 3263),
    @(5025, 18203, 22632, 11244, 0, 1363, 8184),
    @(4582, 20096, 18776, 2485, 0, 1242, 8312),
    @(4077, 21742, 20091, 187, 0, 1009, 8636)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = 2 + $r
    $rowData = $data[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = 13 + $c
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}

# Update the sheet view: scrolled position (topLeftCell I1) and selected cell (T2)
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1
$ws.Range("T2").Select()
